$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G2").Value = 35.73885133333334
$ws.Range("H2").Value = 107.216554
$ws.Range("I2").Value = 0.01949729408921566
$ws.Range("J2").Value = 0.01949729408921566
$ws.Range("K2").Value = 3
$ws.Range("L2").Value = 1
$ws.Range("M2").Value = 7.107333666666666
$ws.Range("N2").Value = 21.322001
$ws.Range("O2").Value = 0.7373665550576455
$ws.Range("P2").Value = 0.7373665550576454
$ws.Range("Q2").Value = 254.0079412893949
$ws.Range("R2").Value = 2286.071471604554
$ws.Range("S2").Value = 0.01437665257551075
$ws.Range("T2").Value = 0.01437665257551074
$ws.Range("G3").Value = 35.73885133333334
$ws.Range("H3").Value = 107.216554
$ws.Range("I3").Value = 0.01949729408921566
$ws.Range("J3").Value = 0.01949729408921566
$ws.Range("O3").Value = 0.1688878844614928
$ws.Range("P3").Value = 0.1688878844614928
$ws.Range("Q3").Value = 58.17847791785356
$ws.Range("R3").Value = 523.606301260682
$ws.Range("S3").Value = 0.003292856751451202
$ws.Range("T3").Value = 0.0032928567514512
$ws.Range("G4").Value = 35.73885133333334
$ws.Range("H4").Value = 107.216554
$ws.Range("I4").Value = 0.01949729408921566
$ws.Range("J4").Value = 0.01949729408921566
$ws.Range("M4").Value = 0.8135026666666666
$ws.Range("N4").Value = 2.440508
$ws.Range("O4").Value = 0.08439869112428164
$ws.Range("P4").Value = 0.08439869112428162
$ws.Range("Q4").Value = 29.07365086327022
$ws.Range("R4").Value = 261.662857769432
$ws.Range("S4").Value = 0.001645546101594995
$ws.Range("T4").Value = 0.001645546101594994
$ws.Range("G5").Value = 35.73885133333334
$ws.Range("H5").Value = 107.216554
$ws.Range("I5").Value = 0.01949729408921566
$ws.Range("J5").Value = 0.01949729408921566
$ws.Range("K5").Value = 2
$ws.Range("L5").Value = 0.6666666666666666
$ws.Range("M5").Value = 0.09009266666666667
$ws.Range("N5").Value = 0.270278
$ws.Range("O5").Value = 0.009346869356580103
$ws.Range("P5").Value = 0.009346869356580103
$ws.Range("Q5").Value = 3.219808420223556
$ws.Range("R5").Value = 28.978275782012
$ws.Range("S5").Value = 0.0001822386606587202
$ws.Range("T5").Value = 0.0001822386606587202
$ws.Range("H6").Value = 5067.86792
$ws.Range("I6").Value = 0.9215900675332435
$ws.Range("J6").Value = 0.9215900675332435
$ws.Range("K6").Value = 3
$ws.Range("L6").Value = 1
$ws.Range("M6").Value = 7.107333666666666
$ws.Range("N6").Value = 21.322001
$ws.Range("O6").Value = 0.7373665550576455
$ws.Range("P6").Value = 0.7373665550576454
$ws.Range("Q6").Value = 12006.34276201199
$ws.Range("R6").Value = 108057.0848581079
$ws.Range("S6").Value = 0.6795496932723306
$ws.Range("T6").Value = 0.6795496932723305
$ws.Range("H7").Value = 5067.86792
$ws.Range("I7").Value = 0.9215900675332435
$ws.Range("J7").Value = 0.9215900675332435
$ws.Range("O7").Value = 0.1688878844614928
$ws.Range("P7").Value = 0.1688878844614928
$ws.Range("S7").Value = 0.1556453968464138
$ws.Range("T7").Value = 0.1556453968464137
$ws.Range("H8").Value = 5067.86792
$ws.Range("I8").Value = 0.9215900675332435
$ws.Range("J8").Value = 0.9215900675332435
$ws.Range("M8").Value = 0.8135026666666666
$ws.Range("N8").Value = 2.440508
$ws.Range("O8").Value = 0.08439869112428164
$ws.Range("P8").Value = 0.08439869112428162
$ws.Range("Q8").Value = 1374.241355744818
$ws.Range("R8").Value = 12368.17220170336
$ws.Range("S8").Value = 0.07778099545294408
$ws.Range("T8").Value = 0.07778099545294406
$ws.Range("H9").Value = 5067.86792
$ws.Range("I9").Value = 0.9215900675332435
$ws.Range("J9").Value = 0.9215900675332435
$ws.Range("K9").Value = 2
$ws.Range("L9").Value = 0.6666666666666666
$ws.Range("M9").Value = 0.09009266666666667
$ws.Range("N9").Value = 0.270278
$ws.Range("O9").Value = 0.009346869356580103
$ws.Range("P9").Value = 0.009346869356580103
$ws.Range("Q9").Value = 152.1925784090844
$ws.Range("R9").Value = 1369.73320568176
$ws.Range("S9").Value = 0.008613981961555061
$ws.Range("T9").Value = 0.008613981961555061
$ws.Range("G10").Value = 93.641553
$ws.Range("H10").Value = 280.924659
$ws.Range("I10").Value = 0.05108605424341119
$ws.Range("J10").Value = 0.05108605424341119
$ws.Range("K10").Value = 3
$ws.Range("L10").Value = 1
$ws.Range("M10").Value = 7.107333666666666
$ws.Range("N10").Value = 21.322001
$ws.Range("O10").Value = 0.7373665550576455
$ws.Range("P10").Value = 0.7373665550576454
$ws.Range("Q10").Value = 665.541762235851
$ws.Range("R10").Value = 5989.875860122659
$ws.Range("S10").Value = 0.03766914782895212
$ws.Range("T10").Value = 0.03766914782895212
$ws.Range("G11").Value = 93.641553
$ws.Range("H11").Value = 280.924659
$ws.Range("I11").Value = 0.05108605424341119
$ws.Range("J11").Value = 0.05108605424341119
$ws.Range("O11").Value = 0.1688878844614928
$ws.Range("P11").Value = 0.1688878844614928
$ws.Range("Q11").Value = 152.436992800683
$ws.Range("R11").Value = 1371.932935206147
$ws.Range("S11").Value = 0.008627815626654784
$ws.Range("T11").Value = 0.008627815626654781
$ws.Range("G12").Value = 93.641553
$ws.Range("H12").Value = 280.924659
$ws.Range("I12").Value = 0.05108605424341119
$ws.Range("J12").Value = 0.05108605424341119
$ws.Range("M12").Value = 0.8135026666666666
$ws.Range("N12").Value = 2.440508
$ws.Range("O12").Value = 0.08439869112428164
$ws.Range("P12").Value = 0.08439869112428162
$ws.Range("Q12").Value = 76.17765307630799
$ws.Range("R12").Value = 685.598877686772
$ws.Range("S12").Value = 0.004311596112847958
$ws.Range("T12").Value = 0.004311596112847958
$ws.Range("G13").Value = 93.641553
$ws.Range("H13").Value = 280.924659
$ws.Range("I13").Value = 0.05108605424341119
$ws.Range("J13").Value = 0.05108605424341119
$ws.Range("K13").Value = 2
$ws.Range("L13").Value = 0.6666666666666666
$ws.Range("M13").Value = 0.09009266666666667
$ws.Range("N13").Value = 0.270278
$ws.Range("O13").Value = 0.009346869356580103
$ws.Range("P13").Value = 0.009346869356580103
$ws.Range("Q13").Value = 8.436417220578001
$ws.Range("R13").Value = 75.92775498520201
$ws.Range("S13").Value = 0.000477494674956329
$ws.Range("T13").Value = 0.000477494674956329
$ws.Range("G14").Value = 14.34625366666667
$ws.Range("H14").Value = 43.038761
$ws.Range("I14").Value = 0.007826584134129748
$ws.Range("J14").Value = 0.007826584134129748
$ws.Range("K14").Value = 3
$ws.Range("L14").Value = 1
$ws.Range("M14").Value = 7.107333666666666
$ws.Range("N14").Value = 21.322001
$ws.Range("O14").Value = 0.7373665550576455
$ws.Range("P14").Value = 0.7373665550576454
$ws.Range("Q14").Value = 101.9636116756401
$ws.Range("R14").Value = 917.6725050807611
$ws.Range("S14").Value = 0.005771061380852077
$ws.Range("T14").Value = 0.005771061380852077
$ws.Range("G15").Value = 14.34625366666667
$ws.Range("H15").Value = 43.038761
$ws.Range("I15").Value = 0.007826584134129748
$ws.Range("J15").Value = 0.007826584134129748
$ws.Range("O15").Value = 0.1688878844614928
$ws.Range("P15").Value = 0.1688878844614928
$ws.Range("Q15").Value = 23.35394594430144
$ws.Range("R15").Value = 210.185513498713
$ws.Range("S15").Value = 0.001321815236973058
$ws.Range("T15").Value = 0.001321815236973057
$ws.Range("G16").Value = 14.34625366666667
$ws.Range("H16").Value = 43.038761
$ws.Range("I16").Value = 0.007826584134129748
$ws.Range("J16").Value = 0.007826584134129748
$ws.Range("M16").Value = 0.8135026666666666
$ws.Range("N16").Value = 2.440508
$ws.Range("O16").Value = 0.08439869112428164
$ws.Range("P16").Value = 0.08439869112428162
$ws.Range("Q16").Value = 11.67071561450978
$ws.Range("R16").Value = 105.036440530588
$ws.Range("S16").Value = 0.0006605534568946198
$ws.Range("T16").Value = 0.0006605534568946197
$ws.Range("G17").Value = 14.34625366666667
$ws.Range("H17").Value = 43.038761
$ws.Range("I17").Value = 0.007826584134129748
$ws.Range("J17").Value = 0.007826584134129748
$ws.Range("K17").Value = 2
$ws.Range("L17").Value = 0.6666666666666666
$ws.Range("M17").Value = 0.09009266666666667
$ws.Range("N17").Value = 0.270278
$ws.Range("O17").Value = 0.009346869356580103
$ws.Range("P17").Value = 0.009346869356580103
$ws.Range("Q17").Value = 1.292492249506445
$ws.Range("R17").Value = 11.632430245558
$ws.Range("S17").Value = 0.00007315405940999336
$ws.Range("T17").Value = 0.00007315405940999336
